$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$sh = $s.Shapes.Item(4)
$tr = $sh.TextFrame.TextRange

# Paragraph 4: "... at offset (...)" -> "... at location (...)"
$para4 = $tr.Paragraphs(4,1)
$run4 = $para4.Runs(1,1)
$run4.Text = "With LM TLV, counter also deeper into the test packet at location (Eth 18, IPv6 40, UDP 8, STAMP 44, TLV Type 4, Total = 114 Byte)"

# Paragraph 5, 3rd run: " / headers in offset" -> " / headers in location"
$para5 = $tr.Paragraphs(5,1)
$run5c = $para5.Runs(3,1)
$run5c.Text = " / headers in location"

# Paragraph 13: "Counter at fixed location - offset (...)" -> "Counter at fixed location (...)"
$para13 = $tr.Paragraphs(13,1)
$run13 = $para13.Runs(1,1)
$run13.Text = "Counter at fixed location (Eth 18, IPv6 40, UDP 8, Seq 4, Total = 70 Byte)"
